$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "count" column (C) values for existing rows (materials used count)
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 1
$ws.Range("C11").Value = 1
$ws.Range("C12").Value = 1

# Row 4 (000004) gets a new usage entry in E/F (AC-Kopplung)
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = "AC-Kopplung"

# Row 2 (000002) gets new aggregated material name and counts
$ws.Range("B2").Value = "Batteriemodul"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 4564596899

# Row 3 (000003) gets new aggregated material name and counts
$ws.Range("B3").Value = "BCU"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 4564596899

# Row 4 (000004) gets new aggregated material name and counts
$ws.Range("B4").Value = "WR"
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 4564596899

# Row 1 (000001) gets new aggregated material name, counts, article number and a usage entry
$ws.Range("B1").Value = "BWR"
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 4635687678
$ws.Range("E1").Value = 5
$ws.Range("F1").Value = "AC-Kopplung"

# Leave the cursor where the author left it when saving
$ws.Range("K7").Select() | Out-Null
